# Generate Report for Handback
# The f68e5e0c-b2a4-47e5-8b0f-52d16917e7ba.md file has now been handed back
# (in sync with en-US), so update the status / datetime / error columns for
# that row on the Overview, zh-cn and de-de sheets.

$wb = $excel.ActiveWorkbook

$statusHandedBack = "Handed back: in sync with en-US"

# --- Overview sheet ---------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = $statusHandedBack
$overview.Range("F3").Value = $statusHandedBack

# --- zh-cn sheet --------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = $statusHandedBack
$zhcn.Range("K3").Value = "2016-08-17 22:47:25"
$zhcn.Range("P3").Value = ""
$zhcn.Columns.Item(16).AutoFit()

# --- de-de sheet --------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = $statusHandedBack
$dede.Range("K3").Value = "2016-08-17 22:47:33"
$dede.Range("P3").Value = ""
$dede.Columns.Item(16).AutoFit()
